$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Row 6
$ws1.Cells.Item(6, 1).Value = 5
$ws1.Cells.Item(6, 2).Value = 500
$ws1.Cells.Item(6, 3).Value = "amo hasan"
$ws1.Cells.Item(6, 4).Value = "MM50"
$ws1.Cells.Item(6, 5).Value = "MM250"
$ws1.Cells.Item(6, 6).Value = "alaki"

# Row 7
$ws1.Cells.Item(7, 1).Value = 6
$ws1.Cells.Item(7, 2).Value = "in chie"
$ws1.Cells.Item(7, 3).Value = "amo mamad"
$ws1.Cells.Item(7, 4).Value = "SH80"
$ws1.Cells.Item(7, 5).Value = "SH100"
$ws1.Cells.Item(7, 6).Value = "dolaki"

# Row 8
$ws1.Cells.Item(8, 1).Value = 7
$ws1.Cells.Item(8, 2).Value = 34234
$ws1.Cells.Item(8, 3).Value = "j100 devices"
$ws1.Cells.Item(8, 4).Value = "MS1500"
$ws1.Cells.Item(8, 5).Value = "MS1600"
$ws1.Cells.Item(8, 6).Value = 41092
$ws1.Cells.Item(8, 6).NumberFormat = "mm/dd/yy"

# Make Sheet1 the active sheet / tab, and set the new selection
$ws1.Activate()
$ws1.Range("E8").Select()
